$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 285, shifting the existing
# rows 285-306 down to 287-308 (matches the new dimension A1:R308).
$ws.Rows.Item(285).Insert()
$ws.Rows.Item(285).Insert()

# --- New row 285 ---
$ws.Range("A285").Value = 10
$ws.Range("B285").Value = "Vega Modelo de Temuco"
$ws.Range("C285").Value = "La Araucanía"
$ws.Range("D285").Value = 44714
$ws.Range("E285").Value = 9
$ws.Range("F285").Value = 100112017
$ws.Range("G285").Value = "Apio"
$ws.Range("H285").Value = "Americana (o)"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 150
$ws.Range("K285").Value = 9000
$ws.Range("L285").Value = 10000
$ws.Range("M285").Value = 9567
$ws.Range("N285").Value = "$/docena de matas"
$ws.Range("O285").Value = "Provincia del Elquí"
$ws.Range("P285").Value = 1594
$ws.Range("Q285").Value = 6
$ws.Range("R285").Value = "Hortaliza"

# --- New row 286 ---
$ws.Range("A286").Value = 10
$ws.Range("B286").Value = "Vega Modelo de Temuco"
$ws.Range("C286").Value = "La Araucanía"
$ws.Range("D286").Value = 44714
$ws.Range("E286").Value = 9
$ws.Range("F286").Value = 100112017
$ws.Range("G286").Value = "Apio"
$ws.Range("H286").Value = "Americana (o)"
$ws.Range("I286").Value = "Segunda"
$ws.Range("J286").Value = 55
$ws.Range("K286").Value = 7000
$ws.Range("L286").Value = 7000
$ws.Range("M286").Value = 7000
$ws.Range("N286").Value = "$/docena de matas"
$ws.Range("O286").Value = "Provincia del Elquí"
$ws.Range("P286").Value = 1167
$ws.Range("Q286").Value = 6
$ws.Range("R286").Value = "Hortaliza"
